$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 53) down to the new
# row 54, then fill in the new values for columns A-E.
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)

$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 2.622852459381209
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 2.447176337618551
